$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 249, shifting existing rows 249:264 down to 250:265
$ws.Rows.Item(249).Insert()

# Populate the newly inserted row 249 with the new weekly price record.
# Columns A-C, E-I, N, O, Q, R mirror the (now shifted-down) row 250 values,
# i.e. the values that used to live in the old row 249.
$ws.Cells.Item(249, 1).Value = 4
$ws.Cells.Item(249, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(249, 3).Value = "Los Lagos"
$ws.Cells.Item(249, 4).Value = 44746
$ws.Cells.Item(249, 5).Value = 10
$ws.Cells.Item(249, 6).Value = 100112044
$ws.Cells.Item(249, 7).Value = "Perejil"
$ws.Cells.Item(249, 8).Value = "Sin especificar"
$ws.Cells.Item(249, 9).Value = "Primera"
$ws.Cells.Item(249, 10).Value = 80
$ws.Cells.Item(249, 11).Value = 5500
$ws.Cells.Item(249, 12).Value = 6000
$ws.Cells.Item(249, 13).Value = 5750
$ws.Cells.Item(249, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(249, 15).Value = "Región Metropolitana"
$ws.Cells.Item(249, 16).Value = 1917
$ws.Cells.Item(249, 17).Value = 3
$ws.Cells.Item(249, 18).Value = "Hortaliza"

# Match the date cell's number format (datetime) used by the rest of column D
$ws.Cells.Item(249, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
